$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($range, $text)
    $escaped = $text.Replace('"', '""')
    $range.Formula = '="' + $escaped + '"'
    $range.Copy()
    $range.PasteSpecial(-4163)
}

# Update ShipmentTracking numbers (column P) for rows 2-22
Set-TextValue $ws.Range("P2") "320018606316"
Set-TextValue $ws.Range("P3") "320018606327"
Set-TextValue $ws.Range("P4") "320018606350"
Set-TextValue $ws.Range("P5") "320018606371"
Set-TextValue $ws.Range("P6") "320018606419"
Set-TextValue $ws.Range("P7") "320018606430"
Set-TextValue $ws.Range("P8") "320018606463"
Set-TextValue $ws.Range("P9") "320018606485"
Set-TextValue $ws.Range("P10") "320018606511"
Set-TextValue $ws.Range("P11") "320018606533"
Set-TextValue $ws.Range("P12") "320018606577"
Set-TextValue $ws.Range("P13") "320018606599"
Set-TextValue $ws.Range("P14") "320018606625"
Set-TextValue $ws.Range("P15") "320018606647"
Set-TextValue $ws.Range("P16") "320018606670"
Set-TextValue $ws.Range("P17") "320018606691"
Set-TextValue $ws.Range("P18") "320018606739"
Set-TextValue $ws.Range("P19") "320018606750"
Set-TextValue $ws.Range("P20") "320018606783"
Set-TextValue $ws.Range("P21") "320018606809"
Set-TextValue $ws.Range("P22") "320018606831"

# Row 20: correct ExpectedRate and Result to match ActualRate
Set-TextValue $ws.Range("Q20") '$62.39'
Set-TextValue $ws.Range("R20") 'PASS'

$excel.CutCopyMode = $false
Write-Output "applied changes"
